# Edit script for Code Metrics.xlsx
# 1) Rename the worksheet from "Sheet1" to "Code Metrics"
# 2) Update metric values across many rows (per diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (this also updates the autofilter defined name reference automatically)
$ws.Name = "Code Metrics"

$ws.Range("G2").Value = 351
$ws.Range("I2").Value = 61
$ws.Range("J2").Value = 1670
$ws.Range("K2").Value = 242
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 385
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 312
$ws.Range("G22").Value = 193
$ws.Range("I22").Value = 35
$ws.Range("J22").Value = 696
$ws.Range("K22").Value = 112
$ws.Range("I23").Value = 5
$ws.Range("J23").Value = 160
$ws.Range("J26").Value = 147
$ws.Range("J27").Value = 71
$ws.Range("J30").Value = 57
$ws.Range("I31").Value = 5
$ws.Range("J31").Value = 79
$ws.Range("J32").Value = 34
$ws.Range("J33").Value = 38
$ws.Range("F34").Value = 54
$ws.Range("G34").Value = 64
$ws.Range("I34").Value = 5
$ws.Range("K34").Value = 19
$ws.Range("F35").Value = 49
$ws.Range("G35").Value = 37
$ws.Range("I35").Value = 3
$ws.Range("K35").Value = 14
$ws.Range("J37").Value = 73
$ws.Range("J44").Value = 39
$ws.Range("F63").Value = 74
$ws.Range("I63").Value = 8
$ws.Range("J65").Value = 32
$ws.Range("J67").Value = 6
$ws.Range("J81").Value = 99
$ws.Range("J90").Value = 14
$ws.Range("J91").Value = 4
$ws.Range("G108").Value = 48
$ws.Range("J108").Value = 319
$ws.Range("K108").Value = 49
$ws.Range("I109").Value = 3
$ws.Range("J109").Value = 18
$ws.Range("F112").Value = 79
$ws.Range("G112").Value = 43
$ws.Range("J112").Value = 288
$ws.Range("K112").Value = 46
$ws.Range("F119").Value = 81
$ws.Range("G119").Value = 8
$ws.Range("I119").Value = 3
$ws.Range("J119").Value = 26
$ws.Range("K119").Value = 8
$ws.Range("F121").Value = 63
$ws.Range("G121").Value = 7
$ws.Range("I121").Value = 3
$ws.Range("J121").Value = 19
$ws.Range("K121").Value = 7
$ws.Range("F137").Value = 93
$ws.Range("G137").Value = 1
$ws.Range("J137").Value = 8
$ws.Range("K137").Value = 1
$ws.Range("I139").Value = 6
$ws.Range("J139").Value = 9
$ws.Range("I144").Value = 26
$ws.Range("J144").Value = 137
$ws.Range("J145").Value = 72
$ws.Range("J148").Value = 31
$ws.Range("I149").Value = 6
$ws.Range("J149").Value = 10
$ws.Range("I151").Value = 5
$ws.Range("J151").Value = 12
